# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45206 to serial date 45208 (2023-10-07 -> 2023-10-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45206) {
        $cell.Value2 = 45208
    }
}
